$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update phone number and rename first_name/other_name to "Shalin"
$ws.Range("A2").Value = 919967488889
$ws.Range("C2").Value = "Shalin"
$ws.Range("D2").Value = "Shalin"

# Row 3: remove the second birthday entry entirely, leaving only the
# (already date-formatted) B3 cell present but empty
$ws.Range("A3").Clear()
$ws.Range("C3").Clear()
$ws.Range("D3").Clear()
$ws.Range("B3").ClearContents()
